$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "15-11-2021 12:27"
$ws.Range("B5").Value = "https://carreras.uleam.edu.ec/citic2021/"
$ws.Range("A6").Value = "15-11-2021 12:39"
$ws.Range("B6").Value = "🙂"
